# Fetch Cell Data based on Rows and Columns Scenario Updated
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after "TestData" and rename it
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TestDataSet"

$xlPasteFormats = -4122

# ---- Formatting: copy header style (row1 of TestData) and body style (row2 of TestData) ----
$ws1.Range("A1:D1").Copy()
$ws2.Range("A2:D2").PasteSpecial($xlPasteFormats)
$ws1.Range("A1").Copy()
$ws2.Range("E2").PasteSpecial($xlPasteFormats)

$ws1.Range("A2:D2").Copy()
$ws2.Range("A3:D5").PasteSpecial($xlPasteFormats)
$ws1.Range("A2").Copy()
$ws2.Range("E3:E5").PasteSpecial($xlPasteFormats)

$ws1.Range("A1:C1").Copy()
$ws2.Range("A8:C8").PasteSpecial($xlPasteFormats)
$ws1.Range("A2:C2").Copy()
$ws2.Range("A9:C10").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---- Table 1: AddCustomerTest header + body (rows 2-5), filled left-to-right, top-to-bottom ----
$ws2.Range("A2").Value = "FirstName"
$ws2.Range("B2").Value = "LastName"
$ws2.Range("C2").Value = "PostCode"
$ws2.Range("D2").Value = "SuccessMessage"
$ws2.Range("E2").Value = "RunMode"

$ws2.Range("A3").Value = "Jack"
$ws2.Range("B3").Value = "Daniel"
$ws2.Range("C3").Value = "JD12345"
$ws2.Range("D3").Value = "Customer added successfully"
$ws2.Range("E3").Value = "Y"

$ws2.Range("A4").Value = "Will"
$ws2.Range("B4").Value = "Smith"
$ws2.Range("C4").Value = "WS12345"
$ws2.Range("D4").Value = "Customer added successfully"
$ws2.Range("E4").Value = "Y"

$ws2.Range("A5").Value = "Blue"
$ws2.Range("B5").Value = "Eye"
$ws2.Range("C5").Value = "BE98765"
$ws2.Range("D5").Value = "Customer added successfully"
$ws2.Range("E5").Value = "Y"

# ---- Table 2: OpenAccountTest header + body (rows 8-10) ----
$ws2.Range("A8").Value = "Customer"
$ws2.Range("B8").Value = "Currency"
$ws2.Range("C8").Value = "SuccessMessage"

$ws2.Range("A9").Value = "Harry Potter"
$ws2.Range("B9").Value = "Rupee"
$ws2.Range("C9").Value = "Account created successfully"

$ws2.Range("A10").Value = "Hermoine Granger"
$ws2.Range("B10").Value = "Dollar"
$ws2.Range("C10").Value = "Account created successfully"

# ---- Table titles, added last ----
$ws2.Range("A1").Value = "AddCustomerTest"
$ws2.Range("A7").Value = "OpenAccountTest"

# ---- Column widths (chosen so the engine's internal width quantization
#      lands as close as possible to the target stored widths of
#      17.5703125 / 9.7109375 / 26.85546875 / 27.140625) ----
$ws2.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 8.833333333333334
$ws2.Columns.Item(3).ColumnWidth = 26.0
$ws2.Columns.Item(4).ColumnWidth = 26.333333333333332

# ---- Selection ----
$ws2.Range("A7").Select()
$ws2.Activate()
